# Refresh cryptocurrency snapshot: prices, 1h volume deltas, and the two
# coins whose rank order flipped since the last run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.621.11"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.440.49"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "492.33"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.76"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +19.57%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.473.17"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.30"
$ws.Range("E10").Value = "  +10.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "2.872.45"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "57.576.07"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.88"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "2.466.16"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.18"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.98"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.37"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "2.557.22"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "0.0₃0811"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.97"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.40"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.35"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.56"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.101"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.56"
$ws.Range("E43").Value = "  +5.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.27"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.68"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.06"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").Value = "1.895.49"
$ws.Range("E51").Value = "  +3.87%  "
